$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.431.86"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.802.58"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.574"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.10%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.96%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0690"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0964"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "1.798.08"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.644"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.50%  "
$ws.Range("D17").Value = "34.413.41"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "0.0₃0790"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("E24").Value = "  +4.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0528"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("D35").Value = "1.395.74"
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.672"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.01%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.959"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "82.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.59%  "
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("E47").Value = "  -4.29%  "
$ws.Range("D48").Value = "1.964.41"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  -1.93%  "
